$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 6704
$ws1.Range("F9").Value  = 6217
$ws1.Range("F12").Value = 1256
$ws1.Range("F18").Value = 364
$ws1.Range("F21").Value = 4535
$ws1.Range("F22").Value = 57
$ws1.Range("F23").Value = 34
$ws1.Range("F25").Value = 62

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 44

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 6704
$ws4.Range("F9").Value  = 6217
$ws4.Range("F12").Value = 1256
$ws4.Range("F18").Value = 364
$ws4.Range("F21").Value = 4535
$ws4.Range("F23").Value = 57
$ws4.Range("F24").Value = 34
$ws4.Range("F26").Value = 62
